$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Format precision data" / "Insert precision data into excel sheet":
# the "Xeon W-10885M" row (row 3) previously held placeholder zeros for
# B:G - replace them with the measured precision values.
$ws.Range("B3").Value = 0.04
$ws.Range("C3").Value = 0.07
$ws.Range("D3").Value = 0.21
$ws.Range("E3").Value = 0.46
$ws.Range("F3").Value = 5.42
$ws.Range("G3").Value = 54.46

# "Update overall data": the "Not enough memory" figure for the largest
# column was rounded to 18.56 - correct it to 18.6.
$ws.Range("G5").Value = 18.6

# The author's last selection before saving moved to A7.
$ws.Range("A7").Select()
